$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.236034750938416
$ws.Range("B1").Value = 2.32744288444519
$ws.Range("C1").Value = 3.320801019668579
$ws.Range("D1").Value = 2.079475402832031
$ws.Range("E1").Value = 1.380895733833313
